$wb = $excel.ActiveWorkbook

# Update the order status on the "Daily Orders" sheet from NEW to COOKING
$ordersWs = $wb.Worksheets.Item("Daily Orders")
$ordersWs.Range("H2").Value = "COOKING"

# Update the corresponding counts on the "Summary" sheet:
# New count goes from 1 to 0, Cooking count goes from 0 to 1
$summaryWs = $wb.Worksheets.Item("Summary")
$summaryWs.Range("B2").Value = 0
$summaryWs.Range("C2").Value = 1
